$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1396.6666
$ws.Range("I70").Value = 1066.6666
$ws.Range("J70").Value = 1506.6666
$ws.Range("K70").Value = 3199.9998
$ws.Range("L70").Value = 4519.9998
$ws.Range("M70").Value = -2929.9998
$ws.Range("N70").Value = -5059.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1396.6666
$ws.Range("I73").Value = 1066.6666
$ws.Range("J73").Value = 1506.6666
$ws.Range("K73").Value = 3199.9998
$ws.Range("L73").Value = 4519.9998
$ws.Range("M73").Value = -2263.9998
$ws.Range("N73").Value = -6391.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 59261444
$ws.Range("I86").Value = 64002080
$ws.Range("K86").Value = 64002080
$ws.Range("M86").Value = -64000957

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 59261444
$ws.Range("I89").Value = 64002080
$ws.Range("K89").Value = 320010400
$ws.Range("M89").Value = -320004784

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5768.25
$ws.Range("J112").Value = 6126.1333
$ws.Range("L112").Value = 18378.3999
$ws.Range("N112").Value = -20594.3999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3152.152
$ws.Range("I137").Value = 3354.3928
$ws.Range("K137").Value = 10063.1784
$ws.Range("M137").Value = -7513.178400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2601.822
$ws.Range("I138").Value = 1467.6052
$ws.Range("J138").Value = 3833.257
$ws.Range("K138").Value = 4402.8156
$ws.Range("L138").Value = 11499.771
$ws.Range("M138").Value = 737.1844000000001
$ws.Range("N138").Value = -21779.771

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 564446.1
$ws.Range("I32").Value = 625458.9399999999
$ws.Range("J32").Value = 24046.857
$ws.Range("K32").Value = 625458.9399999999
$ws.Range("L32").Value = 24046.857
$ws.Range("M32").Value = -625171.9399999999
$ws.Range("N32").Value = -24620.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2802.7551
$ws.Range("I61").Value = 2348.6155
$ws.Range("J61").Value = 3316.1304
$ws.Range("K61").Value = 2348.6155
$ws.Range("L61").Value = 3316.1304
$ws.Range("M61").Value = -2136.6155
$ws.Range("N61").Value = -3740.1304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2480.84
$ws.Range("I74").Value = 2316.0952
$ws.Range("J74").Value = 3345.75
$ws.Range("K74").Value = 2316.0952
$ws.Range("L74").Value = 3345.75
$ws.Range("M74").Value = -1442.0952
$ws.Range("N74").Value = -5093.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2480.84
$ws.Range("I77").Value = 2316.0952
$ws.Range("J77").Value = 3345.75
$ws.Range("K77").Value = 11580.476
$ws.Range("L77").Value = 16728.75
$ws.Range("M77").Value = -7212.476000000001
$ws.Range("N77").Value = -25464.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2689.5083
$ws.Range("I132").Value = 1758.8445
$ws.Range("J132").Value = 5307
$ws.Range("K132").Value = 5276.5335
$ws.Range("L132").Value = 15921
$ws.Range("M132").Value = -2746.5335
$ws.Range("N132").Value = -20981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2802.7551
$ws.Range("I136").Value = 2348.6155
$ws.Range("J136").Value = 3316.1304
$ws.Range("K136").Value = 7045.8465
$ws.Range("L136").Value = 9948.3912
$ws.Range("M136").Value = -4495.8465
$ws.Range("N136").Value = -15048.3912

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1991.6666
$ws.Range("I86").Value = 1491.5
$ws.Range("K86").Value = 1491.5
$ws.Range("M86").Value = -368.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1991.6666
$ws.Range("I89").Value = 1491.5
$ws.Range("K89").Value = 7457.5
$ws.Range("M89").Value = -1841.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10762.893
$ws.Range("I31").Value = 2172.4285
$ws.Range("J31").Value = 13626.381
$ws.Range("K31").Value = 2172.4285
$ws.Range("L31").Value = 13626.381
$ws.Range("M31").Value = -1877.4285
$ws.Range("N31").Value = -14216.381

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10762.893
$ws.Range("I34").Value = 2172.4285
$ws.Range("J34").Value = 13626.381
$ws.Range("K34").Value = 2172.4285
$ws.Range("L34").Value = 13626.381
$ws.Range("M34").Value = -1970.4285
$ws.Range("N34").Value = -14030.381

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1191
$ws.Range("I105").Value = 1036.5
$ws.Range("K105").Value = 1036.5
$ws.Range("M105").Value = 710.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2083808
$ws.Range("I107").Value = 3906473.5
$ws.Range("K107").Value = 3906473.5
$ws.Range("M107").Value = -3904553.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 18818304
$ws.Range("I132").Value = 27778866
$ws.Range("J132").Value = 6411370
$ws.Range("K132").Value = 83336598
$ws.Range("L132").Value = 19234110
$ws.Range("M132").Value = -83334068
$ws.Range("N132").Value = -19239170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1858.3182
$ws.Range("I5").Value = 609.75
$ws.Range("J5").Value = 2571.7856
$ws.Range("K5").Value = 1829.25
$ws.Range("L5").Value = 7715.3568
$ws.Range("M5").Value = -1717.25
$ws.Range("N5").Value = -7939.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 3533.3333
$ws.Range("J96").Value = 3533.3333
$ws.Range("L96").Value = 10599.9999
$ws.Range("N96").Value = -14717.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3064.5366
$ws.Range("J122").Value = 6179.5264
$ws.Range("L122").Value = 55615.7376
$ws.Range("N122").Value = -60515.7376

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5273.607
$ws.Range("J131").Value = 6315.0435
$ws.Range("L131").Value = 18945.1305
$ws.Range("N131").Value = -29025.1305

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1858.3182
$ws.Range("I135").Value = 609.75
$ws.Range("J135").Value = 2571.7856
$ws.Range("K135").Value = 5487.75
$ws.Range("L135").Value = 23146.0704
$ws.Range("M135").Value = -2952.75
$ws.Range("N135").Value = -28216.0704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 325153
$ws.Range("J51").Value = 325153
$ws.Range("L51").Value = 325153
$ws.Range("N51").Value = -326171

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4330482.5
$ws.Range("I16").Value = 1445.5
$ws.Range("K16").Value = 1445.5
$ws.Range("M16").Value = -1275.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 60083
$ws.Range("J63").Value = 60083
$ws.Range("L63").Value = 60083
$ws.Range("N63").Value = -61331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 60083
$ws.Range("J66").Value = 60083
$ws.Range("L66").Value = 180249
$ws.Range("N66").Value = -186489

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 679.3
$ws.Range("I107").Value = 615.8889
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 1847.6667
$ws.Range("L107").Value = 3750
$ws.Range("M107").Value = 72.33329999999978
$ws.Range("N107").Value = -7590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2605837
$ws.Range("I132").Value = 1652.5853
$ws.Range("J132").Value = 7248078.5
$ws.Range("K132").Value = 4957.7559
$ws.Range("L132").Value = 21744235.5
$ws.Range("M132").Value = -2427.7559
$ws.Range("N132").Value = -21749295.5
